$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix row 10 (Objetivos:) content - was incorrectly showing the docente name,
#    now shows the actual Portuguese objectives text.
$ws.Range('B10').Value2 = 'Ter uma maior conscientização entre os alunos sobre questões ligadas à área de sustentabilidadeCompreender o papel da engenharia e da tecnologia no desenvolvimento sustentável;Conhecer os métodos, ferramentas e incentivos para o desenvolvimento sustentável do sistema de produtos-serviçosEstabelecer uma compreensão clara do papel e do impacto de vários aspectos das decisões de engenharia sobre problemas ambientais, sociais e econômicos.'
$ws.Range('C10').Value2 = 'Ter uma maior conscientização entre os alunos sobre questões ligadas à área de sustentabilidadeCompreender o papel da engenharia e da tecnologia no desenvolvimento sustentável;Conhecer os métodos, ferramentas e incentivos para o desenvolvimento sustentável do sistema de produtos-serviçosEstabelecer uma compreensão clara do papel e do impacto de vários aspectos das decisões de engenharia sobre problemas ambientais, sociais e econômicos.'

# 2. Insert a new row at 13 (after "Docentes responsaveis:") to hold the
#    docente identification that was misplaced under "Objetivos:".
$ws.Rows('13').Insert()
$ws.Range('A13').Clear()
$ws.Range('B14').Copy()
$ws.Range('B13').PasteSpecial(-4122)
$ws.Range('C14').Copy()
$ws.Range('C13').PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range('B13').Value2 = '5840535 - Messias Borges Silva'
$ws.Range('C13').Value2 = '5840535 - Messias Borges Silva'

# 3. Row 14 (now "Programa resumido:") gets the new Portuguese short-syllabus text.
$ws.Range('B14').Value2 = 'Sustentabilidade. Protocolos ambientais. Questões ambientais. Recursos naturais e sua poluição, créditos de carbono, conceito de resíduos zero, ISO 14000, análise do ciclo de vida, estudos de avaliação de impacto ambiental, habitat sustentável, Fontes convencionais e renováveis, Tecnologia e desenvolvimento sustentável, Urbanização sustentável, Ecologia Industrial.'
$ws.Range('C14').Value2 = 'Sustentabilidade. Protocolos ambientais. Questões ambientais. Recursos naturais e sua poluição, créditos de carbono, conceito de resíduos zero, ISO 14000, análise do ciclo de vida, estudos de avaliação de impacto ambiental, habitat sustentável, Fontes convencionais e renováveis, Tecnologia e desenvolvimento sustentável, Urbanização sustentável, Ecologia Industrial.'

# 4. Row 16 (now "Programa:") gets the new, full Portuguese syllabus text
#    (previously held a stray date value).
$ws.Range('B16').Value2 = 'Sustentabilidade – necessidade, conceito, desafios, Protocolos ambientais,Questões ambientais globais, regionais e locais, Recursos naturais e sua poluição, Créditos de carbono, Conceito de resíduos zero, ISO 14000,Análise do ciclo de vida, estudos de avaliação de impacto ambiental, habitat sustentável, , Materiais verdes, Energia, Fontes convencionais e renováveis, Tecnologia e desenvolvimento sustentável, Urbanização sustentável, Ecologia Industrial'
$ws.Range('C16').Value2 = 'Sustentabilidade – necessidade, conceito, desafios, Protocolos ambientais,Questões ambientais globais, regionais e locais, Recursos naturais e sua poluição, Créditos de carbono, Conceito de resíduos zero, ISO 14000,Análise do ciclo de vida, estudos de avaliação de impacto ambiental, habitat sustentável, , Materiais verdes, Energia, Fontes convencionais e renováveis, Tecnologia e desenvolvimento sustentável, Urbanização sustentável, Ecologia Industrial'

# 5. Realign the evaluation block content that had drifted one row from its label:
#    row 19 "Metodo:" now correctly shows the teaching method text.
$ws.Range('B19').Value2 = 'Aulas Expositivas; trabalhos e seminários.'
$ws.Range('C19').Value2 = 'Aulas Expositivas; trabalhos e seminários.'

#    row 20 "Criterio:" now correctly shows the assessment criteria text.
$ws.Range('B20').Value2 = 'Avaliação dos trabalhos e apresentações ao longo do semestre'
$ws.Range('C20').Value2 = 'Avaliação dos trabalhos e apresentações ao longo do semestre'

#    row 21 "Norma de recuperacao:" now correctly shows the recovery-grade formula.
$ws.Range('B21').Value2 = 'NF = (MF + PR)/ 2 , onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota da recuperação'
$ws.Range('C21').Value2 = 'NF = (MF + PR)/ 2 , onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota da recuperação'

# 6. New row 22 "Bibliografia:" with its reference list.
$ws.Range('B22').Value2 = 'ALLEN, D.T., SHONNARD, D.R. , Sustainable Engineering :concepts, design and case studies, Prentice Hall, 2015BLOKDIJK, G. , ISO14000 - Simple Steps to Win, Insights and Opportunities for Maxing out Success, Complete Publishing, 2015LAVE, L.B., HENDRICKSON, C.T. , Environmental Life Cycle Assessment of Goods and Services, Ed John Hopkins, 2006'
$ws.Range('C22').Value2 = 'ALLEN, D.T., SHONNARD, D.R. , Sustainable Engineering :concepts, design and case studies, Prentice Hall, 2015BLOKDIJK, G. , ISO14000 - Simple Steps to Win, Insights and Opportunities for Maxing out Success, Complete Publishing, 2015LAVE, L.B., HENDRICKSON, C.T. , Environmental Life Cycle Assessment of Goods and Services, Ed John Hopkins, 2006'

